$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = [double]"11.37596266666667"
$ws.Range("H2").Value = [double]"34.127888"
$ws.Range("I2").Value = [double]"0.05604480707695051"
$ws.Range("J2").Value = [double]"0.05604480707695052"
$ws.Range("M2").Value = [double]"1.701929666666667"
$ws.Range("N2").Value = [double]"5.105789"
$ws.Range("O2").Value = [double]"0.02105622887134972"
$ws.Range("P2").Value = [double]"0.02105622887134972"
$ws.Range("Q2").Value = [double]"19.36108834929244"
$ws.Range("R2").Value = [double]"174.249795143632"
$ws.Range("S2").Value = [double]"0.00118009228486291"
$ws.Range("T2").Value = [double]"0.00118009228486291"

$ws.Range("G3").Value = [double]"11.37596266666667"
$ws.Range("H3").Value = [double]"34.127888"
$ws.Range("I3").Value = [double]"0.05604480707695051"
$ws.Range("J3").Value = [double]"0.05604480707695052"
$ws.Range("O3").Value = [double]"0.7732971809418951"
$ws.Range("P3").Value = [double]"0.7732971809418953"
$ws.Range("Q3").Value = [double]"711.042567591312"
$ws.Range("R3").Value = [double]"6399.383108321808"
$ws.Range("S3").Value = [double]"0.04333929131903821"
$ws.Range("T3").Value = [double]"0.04333929131903821"

$ws.Range("G4").Value = [double]"11.37596266666667"
$ws.Range("H4").Value = [double]"34.127888"
$ws.Range("I4").Value = [double]"0.05604480707695051"
$ws.Range("J4").Value = [double]"0.05604480707695052"
$ws.Range("K4").Value = [double]"3"
$ws.Range("L4").Value = [double]"1"
$ws.Range("M4").Value = [double]"0.4338690000000001"
$ws.Range("N4").Value = [double]"1.301607"
$ws.Range("O4").Value = [double]"0.005367815805265532"
$ws.Range("P4").Value = [double]"0.005367815805265533"
$ws.Range("Q4").Value = [double]"4.935677546224"
$ws.Range("R4").Value = [double]"44.421097916016"
$ws.Range("S4").Value = [double]"0.0003008382012307125"
$ws.Range("T4").Value = [double]"0.0003008382012307126"

$ws.Range("G5").Value = [double]"11.37596266666667"
$ws.Range("H5").Value = [double]"34.127888"
$ws.Range("I5").Value = [double]"0.05604480707695051"
$ws.Range("J5").Value = [double]"0.05604480707695052"
$ws.Range("M5").Value = [double]"15.972384"
$ws.Range("N5").Value = [double]"47.917152"
$ws.Range("O5").Value = [double]"0.1976099128607259"
$ws.Range("P5").Value = [double]"0.1976099128607259"
$ws.Range("Q5").Value = [double]"181.701244081664"
$ws.Range("R5").Value = [double]"1635.311196734976"
$ws.Range("S5").Value = [double]"0.01107500944277238"
$ws.Range("T5").Value = [double]"0.01107500944277239"

$ws.Range("G6").Value = [double]"11.37596266666667"
$ws.Range("H6").Value = [double]"34.127888"
$ws.Range("I6").Value = [double]"0.05604480707695051"
$ws.Range("J6").Value = [double]"0.05604480707695052"
$ws.Range("K6").Value = [double]"3"
$ws.Range("L6").Value = [double]"1"
$ws.Range("M6").Value = [double]"0.2157183333333333"
$ws.Range("N6").Value = [double]"0.647155"
$ws.Range("O6").Value = [double]"0.002668861520763652"
$ws.Range("P6").Value = [double]"0.002668861520763652"
$ws.Range("Q6").Value = [double]"2.454003706515556"
$ws.Range("R6").Value = [double]"22.08603335864"
$ws.Range("S6").Value = [double]"0.0001495758290462957"
$ws.Range("T6").Value = [double]"0.0001495758290462957"

$ws.Range("I7").Value = [double]"0.765548861900355"
$ws.Range("J7").Value = [double]"0.7655488619003551"
$ws.Range("M7").Value = [double]"1.701929666666667"
$ws.Range("N7").Value = [double]"5.105789"
$ws.Range("O7").Value = [double]"0.02105622887134972"
$ws.Range("P7").Value = [double]"0.02105622887134972"
$ws.Range("Q7").Value = [double]"264.4644512845299"
$ws.Range("R7").Value = [double]"2380.180061560769"
$ws.Range("S7").Value = [double]"0.01611957204837517"
$ws.Range("T7").Value = [double]"0.01611957204837517"

$ws.Range("I8").Value = [double]"0.765548861900355"
$ws.Range("J8").Value = [double]"0.7655488619003551"
$ws.Range("O8").Value = [double]"0.7732971809418951"
$ws.Range("P8").Value = [double]"0.7732971809418953"
$ws.Range("S8").Value = [double]"0.5919967767808207"
$ws.Range("T8").Value = [double]"0.5919967767808209"

$ws.Range("I9").Value = [double]"0.765548861900355"
$ws.Range("J9").Value = [double]"0.7655488619003551"
$ws.Range("K9").Value = [double]"3"
$ws.Range("L9").Value = [double]"1"
$ws.Range("M9").Value = [double]"0.4338690000000001"
$ws.Range("N9").Value = [double]"1.301607"
$ws.Range("O9").Value = [double]"0.005367815805265532"
$ws.Range("P9").Value = [double]"0.005367815805265533"
$ws.Range("Q9").Value = [double]"67.41931189148301"
$ws.Range("R9").Value = [double]"606.7738070233471"
$ws.Range("S9").Value = [double]"0.004109325280611765"
$ws.Range("T9").Value = [double]"0.004109325280611767"

$ws.Range("I10").Value = [double]"0.765548861900355"
$ws.Range("J10").Value = [double]"0.7655488619003551"
$ws.Range("M10").Value = [double]"15.972384"
$ws.Range("N10").Value = [double]"47.917152"
$ws.Range("O10").Value = [double]"0.1976099128607259"
$ws.Range("P10").Value = [double]"0.1976099128607259"
$ws.Range("Q10").Value = [double]"2481.963769125088"
$ws.Range("R10").Value = [double]"22337.67392212579"
$ws.Range("S10").Value = [double]"0.151280043890757"
$ws.Range("T10").Value = [double]"0.1512800438907571"

$ws.Range("I11").Value = [double]"0.765548861900355"
$ws.Range("J11").Value = [double]"0.7655488619003551"
$ws.Range("K11").Value = [double]"3"
$ws.Range("L11").Value = [double]"1"
$ws.Range("M11").Value = [double]"0.2157183333333333"
$ws.Range("N11").Value = [double]"0.647155"
$ws.Range("O11").Value = [double]"0.002668861520763652"
$ws.Range("P11").Value = [double]"0.002668861520763652"
$ws.Range("Q11").Value = [double]"33.52067466380611"
$ws.Range("R11").Value = [double]"301.686071974255"
$ws.Range("S11").Value = [double]"0.002043143899790265"
$ws.Range("T11").Value = [double]"0.002043143899790265"

$ws.Range("G12").Value = [double]"11.89345866666667"
$ws.Range("H12").Value = [double]"35.680376"
$ws.Range("I12").Value = [double]"0.05859430238850571"
$ws.Range("J12").Value = [double]"0.05859430238850571"
$ws.Range("M12").Value = [double]"1.701929666666667"
$ws.Range("N12").Value = [double]"5.105789"
$ws.Range("O12").Value = [double]"0.02105622887134972"
$ws.Range("P12").Value = [double]"0.02105622887134972"
$ws.Range("Q12").Value = [double]"20.24183014407378"
$ws.Range("R12").Value = [double]"182.176471296664"
$ws.Range("S12").Value = [double]"0.00123377504164945"
$ws.Range("T12").Value = [double]"0.00123377504164945"

$ws.Range("G13").Value = [double]"11.89345866666667"
$ws.Range("H13").Value = [double]"35.680376"
$ws.Range("I13").Value = [double]"0.05859430238850571"
$ws.Range("J13").Value = [double]"0.05859430238850571"
$ws.Range("O13").Value = [double]"0.7732971809418951"
$ws.Range("P13").Value = [double]"0.7732971809418953"
$ws.Range("Q13").Value = [double]"743.3881101480241"
$ws.Range("R13").Value = [double]"6690.492991332217"
$ws.Range("S13").Value = [double]"0.04531080885628842"
$ws.Range("T13").Value = [double]"0.04531080885628843"

$ws.Range("G14").Value = [double]"11.89345866666667"
$ws.Range("H14").Value = [double]"35.680376"
$ws.Range("I14").Value = [double]"0.05859430238850571"
$ws.Range("J14").Value = [double]"0.05859430238850571"
$ws.Range("K14").Value = [double]"3"
$ws.Range("L14").Value = [double]"1"
$ws.Range("M14").Value = [double]"0.4338690000000001"
$ws.Range("N14").Value = [double]"1.301607"
$ws.Range("O14").Value = [double]"0.005367815805265532"
$ws.Range("P14").Value = [double]"0.005367815805265533"
$ws.Range("Q14").Value = [double]"5.160203018248001"
$ws.Range("R14").Value = [double]"46.44182716423201"
$ws.Range("S14").Value = [double]"0.0003145234224595288"
$ws.Range("T14").Value = [double]"0.0003145234224595289"

$ws.Range("G15").Value = [double]"11.89345866666667"
$ws.Range("H15").Value = [double]"35.680376"
$ws.Range("I15").Value = [double]"0.05859430238850571"
$ws.Range("J15").Value = [double]"0.05859430238850571"
$ws.Range("M15").Value = [double]"15.972384"
$ws.Range("N15").Value = [double]"47.917152"
$ws.Range("O15").Value = [double]"0.1976099128607259"
$ws.Range("P15").Value = [double]"0.1976099128607259"
$ws.Range("Q15").Value = [double]"189.966888912128"
$ws.Range("R15").Value = [double]"1709.702000209152"
$ws.Range("S15").Value = [double]"0.01157881498912764"
$ws.Range("T15").Value = [double]"0.01157881498912764"

$ws.Range("G16").Value = [double]"11.89345866666667"
$ws.Range("H16").Value = [double]"35.680376"
$ws.Range("I16").Value = [double]"0.05859430238850571"
$ws.Range("J16").Value = [double]"0.05859430238850571"
$ws.Range("K16").Value = [double]"3"
$ws.Range("L16").Value = [double]"1"
$ws.Range("M16").Value = [double]"0.2157183333333333"
$ws.Range("N16").Value = [double]"0.647155"
$ws.Range("O16").Value = [double]"0.002668861520763652"
$ws.Range("P16").Value = [double]"0.002668861520763652"
$ws.Range("Q16").Value = [double]"2.565637081142222"
$ws.Range("R16").Value = [double]"23.09073373028"
$ws.Range("S16").Value = [double]"0.0001563800789806727"
$ws.Range("T16").Value = [double]"0.0001563800789806727"

$ws.Range("G17").Value = [double]"23.69116533333333"
$ws.Range("H17").Value = [double]"71.07349600000001"
$ws.Range("I17").Value = [double]"0.1167168730630039"
$ws.Range("J17").Value = [double]"0.1167168730630039"
$ws.Range("M17").Value = [double]"1.701929666666667"
$ws.Range("N17").Value = [double]"5.105789"
$ws.Range("O17").Value = [double]"0.02105622887134972"
$ws.Range("P17").Value = [double]"0.02105622887134972"
$ws.Range("Q17").Value = [double]"40.32069711870489"
$ws.Range("R17").Value = [double]"362.886274068344"
$ws.Range("S17").Value = [double]"0.002457617192362883"
$ws.Range("T17").Value = [double]"0.002457617192362883"

$ws.Range("G18").Value = [double]"23.69116533333333"
$ws.Range("H18").Value = [double]"71.07349600000001"
$ws.Range("I18").Value = [double]"0.1167168730630039"
$ws.Range("J18").Value = [double]"0.1167168730630039"
$ws.Range("O18").Value = [double]"0.7732971809418951"
$ws.Range("P18").Value = [double]"0.7732971809418953"
$ws.Range("Q18").Value = [double]"1480.791342362904"
$ws.Range("R18").Value = [double]"13327.12208126614"
$ws.Range("S18").Value = [double]"0.09025682890797393"
$ws.Range("T18").Value = [double]"0.09025682890797396"

$ws.Range("G19").Value = [double]"23.69116533333333"
$ws.Range("H19").Value = [double]"71.07349600000001"
$ws.Range("I19").Value = [double]"0.1167168730630039"
$ws.Range("J19").Value = [double]"0.1167168730630039"
$ws.Range("K19").Value = [double]"3"
$ws.Range("L19").Value = [double]"1"
$ws.Range("M19").Value = [double]"0.4338690000000001"
$ws.Range("N19").Value = [double]"1.301607"
$ws.Range("O19").Value = [double]"0.005367815805265532"
$ws.Range("P19").Value = [double]"0.005367815805265533"
$ws.Range("Q19").Value = [double]"10.278862212008"
$ws.Range("R19").Value = [double]"92.50975990807203"
$ws.Range("S19").Value = [double]"0.0006265146759687631"
$ws.Range("T19").Value = [double]"0.0006265146759687634"

$ws.Range("G20").Value = [double]"23.69116533333333"
$ws.Range("H20").Value = [double]"71.07349600000001"
$ws.Range("I20").Value = [double]"0.1167168730630039"
$ws.Range("J20").Value = [double]"0.1167168730630039"
$ws.Range("M20").Value = [double]"15.972384"
$ws.Range("N20").Value = [double]"47.917152"
$ws.Range("O20").Value = [double]"0.1976099128607259"
$ws.Range("P20").Value = [double]"0.1976099128607259"
$ws.Range("Q20").Value = [double]"378.404390111488"
$ws.Range("R20").Value = [double]"3405.639511003392"
$ws.Range("S20").Value = [double]"0.02306441111535661"
$ws.Range("T20").Value = [double]"0.02306441111535661"

$ws.Range("G21").Value = [double]"23.69116533333333"
$ws.Range("H21").Value = [double]"71.07349600000001"
$ws.Range("I21").Value = [double]"0.1167168730630039"
$ws.Range("J21").Value = [double]"0.1167168730630039"
$ws.Range("K21").Value = [double]"3"
$ws.Range("L21").Value = [double]"1"
$ws.Range("M21").Value = [double]"0.2157183333333333"
$ws.Range("N21").Value = [double]"0.647155"
$ws.Range("O21").Value = [double]"0.002668861520763652"
$ws.Range("P21").Value = [double]"0.002668861520763652"
$ws.Range("Q21").Value = [double]"5.110618700431112"
$ws.Range("R21").Value = [double]"45.99556830388001"
$ws.Range("S21").Value = [double]"0.0003115011713417068"
$ws.Range("T21").Value = [double]"0.0003115011713417068"

$ws.Range("G22").Value = [double]"0.628254"
$ws.Range("H22").Value = [double]"1.884762"
$ws.Range("I22").Value = [double]"0.003095155571184698"
$ws.Range("J22").Value = [double]"0.003095155571184698"
$ws.Range("M22").Value = [double]"1.701929666666667"
$ws.Range("N22").Value = [double]"5.105789"
$ws.Range("O22").Value = [double]"0.02105622887134972"
$ws.Range("P22").Value = [double]"0.02105622887134972"
$ws.Range("Q22").Value = [double]"1.069244120802"
$ws.Range("R22").Value = [double]"9.623197087217999"
$ws.Range("S22").Value = [double]"6.517230409929816E-05"
$ws.Range("T22").Value = [double]"6.517230409929817E-05"

$ws.Range("G23").Value = [double]"0.628254"
$ws.Range("H23").Value = [double]"1.884762"
$ws.Range("I23").Value = [double]"0.003095155571184698"
$ws.Range("J23").Value = [double]"0.003095155571184698"
$ws.Range("O23").Value = [double]"0.7732971809418951"
$ws.Range("P23").Value = [double]"0.7732971809418953"
$ws.Range("Q23").Value = [double]"39.268354718538"
$ws.Range("R23").Value = [double]"353.415192466842"
$ws.Range("S23").Value = [double]"0.002393475077773728"
$ws.Range("T23").Value = [double]"0.002393475077773729"

$ws.Range("G24").Value = [double]"0.628254"
$ws.Range("H24").Value = [double]"1.884762"
$ws.Range("I24").Value = [double]"0.003095155571184698"
$ws.Range("J24").Value = [double]"0.003095155571184698"
$ws.Range("K24").Value = [double]"3"
$ws.Range("L24").Value = [double]"1"
$ws.Range("M24").Value = [double]"0.4338690000000001"
$ws.Range("N24").Value = [double]"1.301607"
$ws.Range("O24").Value = [double]"0.005367815805265532"
$ws.Range("P24").Value = [double]"0.005367815805265533"
$ws.Range("Q24").Value = [double]"0.272579934726"
$ws.Range("R24").Value = [double]"2.453219412534"
$ws.Range("S24").Value = [double]"1.661422499476089E-05"
$ws.Range("T24").Value = [double]"1.661422499476089E-05"

$ws.Range("G25").Value = [double]"0.628254"
$ws.Range("H25").Value = [double]"1.884762"
$ws.Range("I25").Value = [double]"0.003095155571184698"
$ws.Range("J25").Value = [double]"0.003095155571184698"
$ws.Range("M25").Value = [double]"15.972384"
$ws.Range("N25").Value = [double]"47.917152"
$ws.Range("O25").Value = [double]"0.1976099128607259"
$ws.Range("P25").Value = [double]"0.1976099128607259"
$ws.Range("Q25").Value = [double]"10.034714137536"
$ws.Range("R25").Value = [double]"90.312427237824"
$ws.Range("S25").Value = [double]"0.0006116334227121985"
$ws.Range("T25").Value = [double]"0.0006116334227121986"

$ws.Range("G26").Value = [double]"0.628254"
$ws.Range("H26").Value = [double]"1.884762"
$ws.Range("I26").Value = [double]"0.003095155571184698"
$ws.Range("J26").Value = [double]"0.003095155571184698"
$ws.Range("K26").Value = [double]"3"
$ws.Range("L26").Value = [double]"1"
$ws.Range("M26").Value = [double]"0.2157183333333333"
$ws.Range("N26").Value = [double]"0.647155"
$ws.Range("O26").Value = [double]"0.002668861520763652"
$ws.Range("P26").Value = [double]"0.002668861520763652"
$ws.Range("Q26").Value = [double]"0.13552590579"
$ws.Range("R26").Value = [double]"1.21973315211"
$ws.Range("S26").Value = [double]"8.260541604712083E-06"
$ws.Range("T26").Value = [double]"8.260541604712085E-06"
